# spotipy/songs.xlsx — append newly scraped tracks to the "songs" list and
# tidy up the sheet view (matches the "modified input file and output to
# excel" commit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append the newly scraped songs (rows 7-11) -----------------------
$ws.Range("A7").Value  = "Hold On We’re Going Home – Drake"
$ws.Range("A8").Value  = "I took a pill in ibiza "
$ws.Range("A9").Value  = "Love yourself – justin bieber"
$ws.Range("A10").Value = "side to side ariana grande"
$ws.Range("A11").Value = "truffle butter"

# --- Column A got a touch narrower after the re-export -----------------
$ws.Columns.Item(1).ColumnWidth = 18.333333333333332

# --- Zoom the sheet view way in, selection moves past the last row -----
$excel.ActiveWindow.Zoom = 262
$null = $ws.Range("A12").Select()
